$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply stock-sheet value updates (stock updated by raj time 12:02)
$ws.Range("B18").Value = 43
$ws.Range("C18").Value = 60.5
$ws.Range("E18").Value = 121
$ws.Range("B54").Value = 205
$ws.Range("C54").Value = 1299.5
$ws.Range("E54").Value = 1452.57
$ws.Range("B57").Value = 73
$ws.Range("C57").Value = 394
$ws.Range("E57").Value = 748.6
$ws.Range("B67").Value = 35
$ws.Range("C67").Value = 103
$ws.Range("E67").Value = 53.56
$ws.Range("B69").Value = 74
$ws.Range("C69").Value = 145.5
$ws.Range("E69").Value = 101.85
$ws.Range("B70").Value = 109
$ws.Range("C70").Value = -8
$ws.Range("E70").Value = -7.2
$ws.Range("B71").Value = 108
$ws.Range("C71").Value = 113.5
$ws.Range("E71").Value = 102.15
$ws.Range("B82").Value = 93
$ws.Range("C82").Value = 97.5
$ws.Range("B83").Value = 160
$ws.Range("C83").Value = 112
$ws.Range("E83").Value = 128.80000000000001
$ws.Range("B85").Value = 169
$ws.Range("C85").Value = 4233
$ws.Range("E85").Value = 5926.2
$ws.Range("B87").Value = 85
$ws.Range("C87").Value = 1019
$ws.Range("E87").Value = 1426.6
$ws.Range("B103").Value = 48
$ws.Range("C103").Value = 463
$ws.Range("B153").Value = 108
$ws.Range("C153").Value = 102
$ws.Range("E153").Value = 280.5
$ws.Range("B156").Value = 151
$ws.Range("C156").Value = 130.25
$ws.Range("E156").Value = 442.85
$ws.Range("B193").Value = 56
$ws.Range("C193").Value = 42
$ws.Range("E193").Value = 126
$ws.Range("B209").Value = 80
$ws.Range("C209").Value = 48
$ws.Range("E209").Value = 182.4
$ws.Range("B213").Value = 73
$ws.Range("C213").Value = 4
$ws.Range("E213").Value = 16.47
$ws.Range("B222").Value = 274
$ws.Range("C222").Value = 99.5
$ws.Range("E222").Value = 447.75
$ws.Range("B286").Value = 50
$ws.Range("C286").Value = 25
$ws.Range("E286").Value = 169.97
$ws.Range("B291").Value = 83
$ws.Range("C291").Value = 78
$ws.Range("E291").Value = 482.04
$ws.Range("B298").Value = 41
$ws.Range("C298").Value = 63.5
$ws.Range("E298").Value = 422.28
$ws.Range("B354").Value = 13
$ws.Range("C354").Value = 21.5
$ws.Range("E354").Value = 263.38
$ws.Range("B363").Value = 15
$ws.Range("C363").Value = -1.5
$ws.Range("E363").Value = -18
$ws.Range("B365").Value = 14
$ws.Range("C365").Value = -3
$ws.Range("E365").Value = -33
$ws.Range("B371").Value = 25
$ws.Range("C371").Value = 29.5
$ws.Range("E371").Value = 295
$ws.Range("B396").Value = 47
$ws.Range("C396").Value = 32
$ws.Range("E396").Value = 349.76
$ws.Range("B490").Value = 82
$ws.Range("C490").Value = 360
$ws.Range("E490").Value = 396
$ws.Range("B493").Value = 60
$ws.Range("C493").Value = 159
$ws.Range("E493").Value = 222.6
$ws.Range("B503").Value = 145
$ws.Range("C503").Value = 621.75
$ws.Range("E503").Value = 932.63
$ws.Range("B552").Value = 24
$ws.Range("C552").Value = 4.5
$ws.Range("E552").Value = 21.38
$ws.Range("B587").Value = 48
$ws.Range("C587").Value = 28.5
$ws.Range("E587").Value = 85.5
$ws.Range("B614").Value = 111
$ws.Range("C614").Value = 968
$ws.Range("E614").Value = 798.6
$ws.Range("B625").Value = 148
$ws.Range("C625").Value = 474.5
$ws.Range("E625").Value = 332.15
$ws.Range("B629").Value = 114
$ws.Range("C629").Value = 1165
$ws.Range("E629").Value = 641.08000000000004
$ws.Range("C634").Value = 62702.42
$ws.Range("E634").Value = 120227.62
